$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.442.38"
$ws.Range("E2").Value = "  -0.80%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.563.40"
$ws.Range("E3").Value = "  -1.17%  "
$ws.Range("E4").Value = "  -0.34%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.28"
$ws.Range("E5").Value = "  +1.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.501"
$ws.Range("E6").Value = "  -0.90%  "
$ws.Range("E7").Value = "  -0.31%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.07"
$ws.Range("E8").Value = "  -0.58%  "
$ws.Range("E9").Value = "  -1.55%  "
$ws.Range("E10").Value = "  +0.07%  "
$ws.Range("E11").Value = "  -0.06%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.784.13"
$ws.Range("E12").Value = "  -1.35%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.575.39"
$ws.Range("E13").Value = "  -0.25%  "
$ws.Range("E14").Value = "  -0.88%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.517"
$ws.Range("E15").Value = "  -2.11%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.49"
$ws.Range("E16").Value = "  +0.57%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.434.51"
$ws.Range("E17").Value = "  -0.72%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "212.77"
$ws.Range("E18").Value = "  -2.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0689"
$ws.Range("E19").Value = "  -0.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.23"
$ws.Range("E20").Value = "  -0.84%  "
$ws.Range("E21").Value = "  -0.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.09"
$ws.Range("E22").Value = "  -0.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.56"
$ws.Range("E23").Value = "  +1.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.04"
$ws.Range("E24").Value = "  +3.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.63"
$ws.Range("E25").Value = "  -0.53%  "
$ws.Range("E26").Value = "  -0.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.66"
$ws.Range("E27").Value = "  -1.18%  "
$ws.Range("E28").Value = "  -0.86%  "
$ws.Range("E29").Value = "  -2.02%  "
$ws.Range("E30").Value = "  -0.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0469"
$ws.Range("E31").Value = "  +1.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.19"
$ws.Range("E32").Value = "  -0.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.373.05"
$ws.Range("E33").Value = "  -0.83%  "
$ws.Range("E34").Value = "  +0.92%  "
$ws.Range("E35").Value = "  +1.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.955"
$ws.Range("E36").Value = "  -1.25%  "
$ws.Range("E37").Value = "  -1.06%  "
$ws.Range("E38").Value = "  +1.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.529"
$ws.Range("E39").Value = "  -1.79%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.820"
$ws.Range("E40").Value = "  +0.26%  "
$ws.Range("E41").Value = "  -0.25%  "
$ws.Range("E42").Value = "  -0.29%  "
$ws.Range("E43").Value = "  +1.98%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.85"
$ws.Range("E44").Value = "  +0.86%  "
$ws.Range("E45").Value = "  -0.57%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.24"
$ws.Range("E46").Value = "  +0.73%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.696.86"
$ws.Range("E47").Value = "  -1.33%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.33"
$ws.Range("E48").Value = "  -2.92%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₇0996"
$ws.Range("E49").Value = "  -0.64%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0955"
$ws.Range("E50").Value = "  -1.62%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0495"
$ws.Range("E51").Value = "  -0.53%  "
